# 检查列表2019 - daily checklist update
# Fills in the remaining checkmarks for 2019-07-14, adds the new rows for
# 2019-07-15 (fully completed) and 2019-07-16 (in-progress / just started).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15: 2019-07-14, finish filling in the day's checklist ---
$row15 = @("√","×","×","√","√","√","×","√","√","√","√","√","√","×")
for ($i = 0; $i -lt $row15.Length; $i++) {
    $ws.Cells.Item(15, 2 + $i).Value = $row15[$i]
}

# --- Row 16: 2019-07-15, new day, fully checked off ---
$ws.Cells.Item(16, 1).Value = 20190715
$row16 = @("√","√","×","√","√","√","√","√","√","√","√","√","√","√")
for ($i = 0; $i -lt $row16.Length; $i++) {
    $ws.Cells.Item(16, 2 + $i).Value = $row16[$i]
}

# --- Row 17: 2019-07-16, new day, just started (some columns still blank) ---
$ws.Cells.Item(17, 1).Value = 20190716
$ws.Cells.Item(17, 2).Value = "√"
$ws.Cells.Item(17, 3).Value = "√"
# D17 (column 4) left blank (not yet recorded)
$ws.Cells.Item(17, 5).Value = "√"
$ws.Cells.Item(17, 6).Value = "√"
$ws.Cells.Item(17, 7).Value = "√"
$ws.Cells.Item(17, 8).Value = "√"
$ws.Cells.Item(17, 9).Value = "√"
$ws.Cells.Item(17, 10).Value = "√"
$ws.Cells.Item(17, 11).Value = "√"
$ws.Cells.Item(17, 12).Value = "√"
# M17 (column 13) left blank (not yet recorded)
$ws.Cells.Item(17, 14).Value = "√"
# O17 (column 15) left blank (not yet recorded)

$ws.Range("O17").Select()
